$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, preserving values such as
# "1.00", "0.0000261", and thousand-grouped strings like "63.094.93"
# that would otherwise be auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "63.094.93", "  -5.82%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.305.71", "  -6.52%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.02%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "547.87", "  -2.68%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "169.14", "  -10.33%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.603", "  -3.81%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.15%  ")
    ,@("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.286.29", "  -6.88%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.606", "  -4.48%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.146", "  -5.81%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "53.73", "  -1.82%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000261", "  -5.41%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "8.78", "  -6.23%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.843.83", "  -6.22%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.116", "  -4.57%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.300.46", "  -6.85%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "17.43", "  -6.50%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "63.076.12", "  -5.95%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "11.42", "  -5.64%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.958", "  -4.26%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "397.60", "  -5.78%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.99", "  -3.15%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "4.32", "  +3.08%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "13.10", "  +6.28%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "81.22", "  -5.14%  ")
    ,@("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "10.66", "  -3.56%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.69", "  -7.54%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "8.52", "  -6.61%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "28.72", "  -5.69%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "6.48", "  -2.63%  ")
    ,@("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "581.87", "  -8.53%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.16", "  -4.98%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.104", "  -6.77%  ")
    ,@("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "57.53", "  -4.87%  ")
    ,@("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.148", "  -1.02%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  +0.19%  ")
    ,@("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "35.40", "  -7.78%  ")
    ,@("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "3.44", "  +2.17%  ")
    ,@("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0724", "  -12.34%  ")
    ,@("TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.361", "  -6.96%  ")
    ,@("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "3.106.72", "  -1.03%  ")
    ,@("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  -0.15%  ")
    ,@("ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "2.75", "  -4.16%  ")
    ,@("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "2.45", "  -7.16%  ")
    ,@("ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.17", "  -5.80%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0396", "  -5.36%  ")
    ,@("WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.58", "  -6.31%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.126", "  -4.93%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "131.38", "  -6.37%  ")
    ,@("THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "7.94", "  -7.44%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
